$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel stores the exact
# string instead of coercing to a floating point number.
$textCells = @("D4","D5","D6","D9","D10","D14","D15","D17","D19","D20","D21","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D36","D37","D38","D39","D41","D44","D45","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '29.389.31'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.841.71'
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '239.16'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '0.6268'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = '24.96'
$ws.Range("E9").Value = '  +2.11%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.2892'
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = '1.837.47'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").Value = '0.6744'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '0.00001031'
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D17").Value = '6.209'
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").Value = '29.419.99'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '233.98'
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").Value = '7.293'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '158.61'
$ws.Range("D25").Value = '8.501'
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").Value = '0.1344'
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("D28").Value = '0.07309'
$ws.Range("E28").Value = '  +12.73%  '
$ws.Range("D29").Value = '1.465'
$ws.Range("E29").Value = '  +4.06%  '
$ws.Range("D30").Value = '1.481'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.031'
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '4.038'
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").Value = '1.816'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '2.572'
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").Value = '0.01843'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").Value = '6.906'
$ws.Range("E38").Value = '  +2.07%  '
$ws.Range("D39").Value = '2.815'
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").Value = '1.233.02'
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("D41").Value = '0.9571'
$ws.Range("E41").Value = '  +4.52%  '
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").Value = '1.997.65'
$ws.Range("E43").Value = '  -0.51%  '
$ws.Range("D44").Value = '100.99'
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").Value = '65.41'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("E46").Value = '  +4.65%  '
$ws.Range("D47").Value = '1.714'
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").Value = '6.952'
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  -2.69%  '
$ws.Range("D50").Value = '8.845'
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").Value = '0.3898'
$ws.Range("E51").Value = '  -1.67%  '
